# Auto-generated Excel COM-interop script to apply cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.242.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.244.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.11"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.12"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.572"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.84"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.34%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.21"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.335.02"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.27%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.586.84"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.833"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.985.33"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0966"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.16"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -7.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.62"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.81"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.96"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.96"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.42"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.27%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.10"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.27"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0797"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.20%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.16%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.18%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.54"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.84"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.45"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -7.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0300"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.747.97"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "82.86"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.192"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.10"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.12"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.68"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.63%  "
